$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 306.875
$ws.Range("I4").Value = 207.85715
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 207.85715
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -93.85714999999999
$ws.Range("N4").Value = -1228

# Row 15
$ws.Range("H15").Value = 1141.2646
$ws.Range("I15").Value = 1141.2646
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 3423.7938
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -3254.7938

# Row 19
$ws.Range("H19").Value = 1990.3
$ws.Range("I19").Value = 1400.5
$ws.Range("J19").Value = 2137.75
$ws.Range("K19").Value = 1400.5
$ws.Range("L19").Value = 2137.75
$ws.Range("M19").Value = -1225.5
$ws.Range("N19").Value = -2487.75

# Row 40
$ws.Range("H40").Value = 2899.6785
$ws.Range("I40").Value = 1614.5385
$ws.Range("J40").Value = 4013.4666
$ws.Range("K40").Value = 1614.5385
$ws.Range("L40").Value = 4013.4666
$ws.Range("M40").Value = -1439.5385
$ws.Range("N40").Value = -4363.4666

# Row 58
$ws.Range("H58").Value = 1708.5714
$ws.Range("I58").Value = 1486.6666
$ws.Range("J58").Value = 1875
$ws.Range("K58").Value = 4459.9998
$ws.Range("L58").Value = 5625
$ws.Range("M58").Value = -4309.9998
$ws.Range("N58").Value = -5925

# Row 107
$ws.Range("H107").Value = 545
$ws.Range("I107").Value = 486.42856
$ws.Range("J107").Value = 750
$ws.Range("K107").Value = 486.42856
$ws.Range("L107").Value = 750
$ws.Range("M107").Value = 1433.57144

# Row 116
$ws.Range("H116").Value = 19114.334
$ws.Range("I116").Value = 2421.25
$ws.Range("J116").Value = 52500.5
$ws.Range("K116").Value = 2421.25
$ws.Range("L116").Value = 52500.5
$ws.Range("M116").Value = 1020.75
$ws.Range("N116").Value = -59384.5

# Row 132
$ws.Range("H132").Value = 29241704
$ws.Range("I132").Value = 5051552.5
$ws.Range("J132").Value = 62503164
$ws.Range("K132").Value = 15154657.5
$ws.Range("L132").Value = 187509492
$ws.Range("M132").Value = -15152127.5
$ws.Range("N132").Value = -187514552

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4553.291
$ws.Range("I32").Value = 3048.6
$ws.Range("J32").Value = 19600.2
$ws.Range("K32").Value = 3048.6
$ws.Range("L32").Value = 19600.2
$ws.Range("M32").Value = -2761.6
$ws.Range("N32").Value = -20174.2

# Row 135
$ws.Range("H135").Value = 30000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 30000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 30000
$ws.Range("N135").Value = -40140

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 2364
$ws.Range("I99").Value = 2578.1538
$ws.Range("J99").Value = 1900
$ws.Range("K99").Value = 2578.1538
$ws.Range("L99").Value = 1900
$ws.Range("M99").Value = -1080.1538

$ws = $wb.Worksheets.Item("CRP")
# Row 25
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = ""

# Row 31
$ws.Range("H31").Value = 14287068
$ws.Range("I31").Value = 1085.3158
$ws.Range("J31").Value = 31251674
$ws.Range("K31").Value = 1085.3158
$ws.Range("L31").Value = 31251674
$ws.Range("M31").Value = -790.3158000000001
$ws.Range("N31").Value = -31252264

# Row 34
$ws.Range("H34").Value = 14287068
$ws.Range("I34").Value = 1085.3158
$ws.Range("J34").Value = 31251674
$ws.Range("K34").Value = 1085.3158
$ws.Range("L34").Value = 31251674
$ws.Range("M34").Value = -883.3158000000001
$ws.Range("N34").Value = -31252078

# Row 99
$ws.Range("H99").Value = 5292950
$ws.Range("I99").Value = 2375
$ws.Range("J99").Value = 9525410
$ws.Range("K99").Value = 2375
$ws.Range("L99").Value = 9525410
$ws.Range("M99").Value = -877
$ws.Range("N99").Value = -9528406

# Row 126
$ws.Range("H126").Value = 5292950
$ws.Range("I126").Value = 2375
$ws.Range("J126").Value = 9525410
$ws.Range("K126").Value = 7125
$ws.Range("L126").Value = 28576230
$ws.Range("M126").Value = -4655
$ws.Range("N126").Value = -28581170

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 2110
$ws.Range("I3").Value = 1665
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 4995
$ws.Range("L3").Value = 9000
$ws.Range("M3").Value = -4883
$ws.Range("N3").Value = -9224

# Row 7
$ws.Range("H7").Value = 22222394
$ws.Range("I7").Value = 50
$ws.Range("J7").Value = 25000188
$ws.Range("K7").Value = 150
$ws.Range("L7").Value = 75000564
$ws.Range("M7").Value = -38
$ws.Range("N7").Value = -75000788

# Row 17
$ws.Range("H17").Value = 1618.7894
$ws.Range("I17").Value = 150.7
$ws.Range("J17").Value = 3250
$ws.Range("K17").Value = 452.1
$ws.Range("L17").Value = 9750
$ws.Range("M17").Value = -283.1
$ws.Range("N17").Value = -10088

# Row 51
$ws.Range("H51").Value = 2709.0908
$ws.Range("I51").Value = 433.33334
$ws.Range("J51").Value = 3562.5
$ws.Range("K51").Value = 1300.00002
$ws.Range("L51").Value = 10687.5
$ws.Range("M51").Value = -840.0000199999999
$ws.Range("N51").Value = -11607.5

# Row 119
$ws.Range("H119").Value = 3864.8572
$ws.Range("I119").Value = 4198
$ws.Range("J119").Value = 3032
$ws.Range("K119").Value = 12594
$ws.Range("L119").Value = 9096
$ws.Range("M119").Value = -7756
$ws.Range("N119").Value = -18772

# Row 123
$ws.Range("H123").Value = 6837.857
$ws.Range("I123").Value = 6010
$ws.Range("J123").Value = 7063.636
$ws.Range("K123").Value = 18030
$ws.Range("L123").Value = 21190.908
$ws.Range("M123").Value = -15580
$ws.Range("N123").Value = -26090.908

# Row 132
$ws.Range("H132").Value = 6338370.5
$ws.Range("I132").Value = 3510612
$ws.Range("J132").Value = 10471249
$ws.Range("K132").Value = 31595508
$ws.Range("L132").Value = 94241241
$ws.Range("M132").Value = -31592978
$ws.Range("N132").Value = -94246301

# Row 139
$ws.Range("H139").Value = 1797.7778
$ws.Range("I139").Value = 1377
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 4131
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = 1009

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 1809.3334
$ws.Range("I113").Value = 1584.8182
$ws.Range("J113").Value = 2056.3
$ws.Range("K113").Value = 1584.8182
$ws.Range("L113").Value = 2056.3
$ws.Range("M113").Value = 585.1818000000001
$ws.Range("N113").Value = -6396.3

# Row 132
$ws.Range("H132").Value = 3921.4285
$ws.Range("I132").Value = 4207.846
$ws.Range("J132").Value = 3456
$ws.Range("K132").Value = 12623.538
$ws.Range("L132").Value = 10368
$ws.Range("M132").Value = -10093.538
$ws.Range("N132").Value = -15428

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 940.1667
$ws.Range("I46").Value = 1120.2858
$ws.Range("J46").Value = 866
$ws.Range("K46").Value = 1120.2858
$ws.Range("L46").Value = 866
$ws.Range("M46").Value = -932.2858000000001
$ws.Range("N46").Value = -1242

# Row 94
$ws.Range("H94").Value = 12999.833
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 12999.833
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 12999.833
$ws.Range("N94").Value = -14351.833

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 493.77274
$ws.Range("I113").Value = 368.92856
$ws.Range("J113").Value = 712.25
$ws.Range("K113").Value = 1106.78568
$ws.Range("L113").Value = 2136.75
$ws.Range("M113").Value = 1063.21432
$ws.Range("N113").Value = -6476.75

# Row 122
$ws.Range("H122").Value = 57436.59
$ws.Range("I122").Value = 4170
$ws.Range("J122").Value = 101825.414
$ws.Range("K122").Value = 12510
$ws.Range("L122").Value = 305476.242
$ws.Range("M122").Value = -10060
$ws.Range("N122").Value = -310376.242
